$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StartSceneConfig")

# New row 13: an extra "LoginCenter" scene, formatted like the "Account"
# row above it (row 12). Column H (OuterPort) is intentionally left blank
# for this scene, so only C:G are touched.
$ws.Range("C12:G12").Copy($ws.Range("C13:G13"))

$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "LoginCenter"
$ws.Range("G13").Value = "LoginCenter"

$ws.Range("F16").Select()
